# Fruta / hortaliza, semanal
# Insert a new data row (weekly update) above the existing row 58, shifting
# all the following rows down by one, and populate the new row with the
# latest market observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 58 (this pushes old rows 58..188 down to 59..189
# and copies the formatting - e.g. the date style on column D - from the row
# that used to be at 58).
$ws.Rows.Item(58).Insert()

# Populate the newly inserted row 58 with the new observation. Columns
# A, B, C, E, F, G, H, I, J, K and T keep the same values as the rest of the
# series for this market/product (Agricola del Norte S.A. de Arica - Manzana).
$ws.Range("A58").Value = 1
$ws.Range("B58").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C58").Value = "Arica y Parinacota"
$ws.Range("D58").Value = 44994
$ws.Range("E58").Value = 15
$ws.Range("F58").Value = "Fruta"
$ws.Range("G58").Value = 100104
$ws.Range("H58").Value = "Frutos de pepita"
$ws.Range("I58").Value = 100104002
$ws.Range("J58").Value = "Manzana"
$ws.Range("K58").Value = "Fuji royal"
$ws.Range("L58").Value = "Segunda"
$ws.Range("M58").Value = 300
$ws.Range("N58").Value = 23000
$ws.Range("O58").Value = 24000
$ws.Range("P58").Value = 23500
$ws.Range("Q58").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R58").Value = "Región Metropolitana"
$ws.Range("S58").Value = 1306
$ws.Range("T58").Value = 18
